$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7's "polite_expressions" (C7) previously held the stray text "nan".
# Clear it back to blank (matching the empty-text pattern already used by
# every other row in that column, C2:C6). A bare "'" forces an empty TEXT
# cell (rather than Excel's default "clear -> blank/number" behaviour), and
# ClearFormats() drops the transient quote-prefix formatting that typing
# "'" would otherwise leave behind.
$ws.Range("C7").Value = "'"
$ws.Range("C7").ClearFormats()

# Append the new annotation as row 8.
$ws.Range("A8").Value = "parisk"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "nan"
$ws.Range("D8").Value = "APC"
$ws.Range("E8").Value = "THE"
$ws.Range("F8").Value = "f5b44bd7-9311-4cfc-b939-3b86c20706ac"
$ws.Range("G8").Value = "SkYXvCR6W_annotated.xlsx"
$ws.Range("H8").Value = "On top of this, I do not enjoy the style the paper is written in, the language is convoluted."
